$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text would otherwise be auto-parsed as a number by Excel;
# format them as Text first so the literal string (e.g. trailing zeros) is preserved.
$textCells = @('D5', 'D8', 'D9', 'D13', 'D15', 'D18', 'D19', 'D21', 'D22', 'D25', 'D26', 'D27', 'D32', 'D33', 'D36', 'D40', 'D43', 'D44', 'D47', 'D48', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '35.701.14'
$ws.Range('E2').Value = '  +3.47%  '

$ws.Range('D3').Value = '1.864.43'
$ws.Range('E3').Value = '  +2.93%  '

$ws.Range('E4').Value = '  +0.40%  '

$ws.Range('D5').Value = '231.48'
$ws.Range('E5').Value = '  +2.60%  '

$ws.Range('E6').Value = '  +3.34%  '

$ws.Range('E7').Value = '  +0.38%  '

$ws.Range('D8').Value = '42.59'
$ws.Range('E8').Value = '  +11.67%  '

$ws.Range('D9').Value = '0.310'
$ws.Range('E9').Value = '  +7.54%  '

$ws.Range('E10').Value = '  +3.18%  '

$ws.Range('E11').Value = '  +4.18%  '

$ws.Range('E12').Value = '  +3.00%  '

$ws.Range('D13').Value = '11.68'
$ws.Range('E13').Value = '  +3.76%  '

$ws.Range('D14').Value = '1.858.28'
$ws.Range('E14').Value = '  +2.69%  '

$ws.Range('D15').Value = '0.683'
$ws.Range('E15').Value = '  +7.91%  '

$ws.Range('E16').Value = '  +7.24%  '

$ws.Range('D17').Value = '35.713.32'
$ws.Range('E17').Value = '  +3.61%  '

$ws.Range('D18').Value = '70.57'
$ws.Range('E18').Value = '  +3.18%  '

$ws.Range('D19').Value = '249.49'
$ws.Range('E19').Value = '  +2.53%  '

$ws.Range('E20').Value = '  +4.14%  '

$ws.Range('D21').Value = '12.34'
$ws.Range('E21').Value = '  +9.97%  '

$ws.Range('D22').Value = '4.76'
$ws.Range('E22').Value = '  +15.32%  '

$ws.Range('E23').Value = '  +0.40%  '

$ws.Range('E24').Value = '  +1.20%  '

$ws.Range('D25').Value = '170.80'
$ws.Range('E25').Value = '  +0.07%  '

$ws.Range('D26').Value = '8.03'
$ws.Range('E26').Value = '  +3.22%  '

$ws.Range('D27').Value = '17.95'
$ws.Range('E27').Value = '  +1.28%  '

$ws.Range('E28').Value = '  +1.75%  '

$ws.Range('E29').Value = '  +16.55%  '

$ws.Range('E30').Value = '  +0.46%  '

$ws.Range('D31').Value = '3.285.92'
$ws.Range('E31').Value = '  +35.24%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '4.10'
$ws.Range('E32').Value = '  +6.11%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.0548'
$ws.Range('E33').Value = '  +6.03%  '

$ws.Range('E34').Value = '  +4.43%  '

$ws.Range('E35').Value = '  +4.46%  '

$ws.Range('D36').Value = '102.51'
$ws.Range('E36').Value = '  +25.15%  '

$ws.Range('E37').Value = '  +7.61%  '

$ws.Range('D38').Value = '1.369.59'

$ws.Range('E39').Value = '  +7.04%  '

$ws.Range('D40').Value = '1.10'
$ws.Range('E40').Value = '  +3.55%  '

$ws.Range('E41').Value = '  +5.08%  '

$ws.Range('E42').Value = '  +6.29%  '

$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '1.27'
$ws.Range('E43').Value = '  +4.26%  '

$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '14.92'
$ws.Range('E44').Value = '  +8.39%  '

$ws.Range('E45').Value = '  +0.81%  '

$ws.Range('E46').Value = '  +0.94%  '

$ws.Range('D47').Value = '6.30'
$ws.Range('E47').Value = '  +8.69%  '

$ws.Range('D48').Value = '0.0521'
$ws.Range('E48').Value = '  +2.45%  '

$ws.Range('D49').Value = '2.032.29'

$ws.Range('D50').Value = '105.02'
$ws.Range('E50').Value = '  +2.36%  '

$ws.Range('E51').Value = '  +0.40%  '
